$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New mini "join" example tables added to the right of the room table ---

# Table 1: F2:H4  (Sid / Dname / Dcount)
$ws.Range("F2").Value = "Sid"
$ws.Range("G2").Value = "Dname"
$ws.Range("H2").Value = "Dcount"
$ws.Range("F3").Value = "s2"
$ws.Range("G3").Value = "EE"
$ws.Range("H3").Value = 140
$ws.Range("F4").Value = "s5"
$ws.Range("G4").Value = "CS"
$ws.Range("H4").Value = 140
$ws.Range("F2:H4").Font.Name = "Fira Code"

# Table 2a: F8:G10 (Sid / Dcount)   Table 2b: I8:J10 (Sid / Dname)
$ws.Range("F8").Value = "Sid"
$ws.Range("G8").Value = "Dcount"
$ws.Range("I8").Value = "Sid"
$ws.Range("J8").Value = "Dname"

$ws.Range("F9").Value = "s2"
$ws.Range("G9").Value = 140
$ws.Range("I9").Value = "s2"
$ws.Range("J9").Value = "EE"

$ws.Range("F10").Value = "s5"
$ws.Range("G10").Value = 140
$ws.Range("I10").Value = "s5"
$ws.Range("J10").Value = "CS"

$ws.Range("F8:G10").Font.Name = "Fira Code"
$ws.Range("I8:J10").Font.Name = "Fira Code"

# Table 3a: F12:G14 (Sid / Dcount)   Table 3b: I12:J14 (Dcount / Dname)
$ws.Range("F12").Value = "Sid"
$ws.Range("G12").Value = "Dcount"
$ws.Range("I12").Value = "Dcount"
$ws.Range("J12").Value = "Dname"

$ws.Range("F13").Value = "s2"
$ws.Range("G13").Value = 140
$ws.Range("I13").Value = 140
$ws.Range("J13").Value = "EE"

$ws.Range("F14").Value = "s5"
$ws.Range("G14").Value = 140
$ws.Range("I14").Value = 140
$ws.Range("J14").Value = "CS"

$ws.Range("F12:G14").Font.Name = "Fira Code"
$ws.Range("I12:J14").Font.Name = "Fira Code"

# Match the selection left behind in the saved file.
$ws.Range("I12:J14").Select()
